# Apply the latest crypto price/volume snapshot to the sheet.
# Numeric-looking "Price" values must be forced to Text format so Excel
# keeps the original decimal-grouped string instead of coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.234.67'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '3.549.78'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '615.31'
$ws.Range('E5').Value = '  +4.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '189.09'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.638'
$ws.Range('E7').Value = '  +2.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.216'
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.96'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000309'
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.76'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').Value = '4.118.02'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '614.39'
$ws.Range('E15').Value = '  +8.04%  '
$ws.Range('E16').Value = '  +3.67%  '
$ws.Range('D17').Value = '70.352.52'
$ws.Range('E17').Value = '  -1.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.21'
$ws.Range('D19').Value = '3.588.09'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.66'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '105.84'
$ws.Range('E23').Value = '  +11.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.71'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.14'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('E26').Value = '  +3.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.00'
$ws.Range('E27').Value = '  -3.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.17'
$ws.Range('E28').Value = '  +10.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.56'
$ws.Range('E29').Value = '  +6.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('E30').Value = '  -2.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.54'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  +1.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.34'
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('E34').Value = '  +12.90%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '539.68'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.17'
$ws.Range('E36').Value = '  -5.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.401'
$ws.Range('E38').Value = '  -4.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.30'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').Value = '0.0₃0782'
$ws.Range('E40').Value = '  -3.87%  '
$ws.Range('E41').Value = '  +3.64%  '
$ws.Range('D42').Value = '3.546.06'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('E43').Value = '  +2.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0467'
$ws.Range('E44').Value = '  +4.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.97'
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.144'
$ws.Range('E46').Value = '  +4.31%  '
$ws.Range('E47').Value = '  -2.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.01'
$ws.Range('E48').Value = '  -5.37%  '
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.68'
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('E51').Value = '  -5.84%  '

Write-Host "Applied all changes"
